# globalcorp_org_chart.xlsx -- "Add files via upload" edit
#
# Employees sheet:
#   - 4 departed employees removed (Hiroshi Tanaka, Luca Romano, Mei Lin,
#     Fatima Ndiaye)
#   - Frederik Jensen (previously no Title/Department) becomes
#     "General Counsel" in "Legal"
#   - remaining employees re-sorted by Title
#   - mailto hyperlinks rebuilt for the employees that still have one
#     (Jensen keeps the Hyperlink look on C but no live link)
#   - columns C/D get explicit widths
#
# Departments sheet:
#   - duplicate "Human Resources" row removed
#
# View state: Employees becomes the active/selected tab.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Employees")
$ws2 = $wb.Worksheets.Item("Departments")

# ---------------------------------------------------------------------
# Employees sheet
# ---------------------------------------------------------------------

# Drop all existing hyperlinks up front -- they get rebuilt after the
# sort below, once everybody is on their final row.
$ws1.Hyperlinks.Delete()

# Delete departed employees bottom-up so earlier row numbers stay valid.
$ws1.Rows.Item(11).Delete()   # Fatima Ndiaye
$ws1.Rows.Item(8).Delete()    # Mei Lin
$ws1.Rows.Item(7).Delete()    # Luca Romano
$ws1.Rows.Item(4).Delete()    # Hiroshi Tanaka

# Frederik Jensen (now on row 7) picks up a real Title/Department.
$ws1.Range("D7").Value = "General Counsel"
$ws1.Range("E7").Value = "Legal"

# Sort the employee block (A2:E10, header in row 1) by Title ascending --
# matches the sortState left behind in the workbook.
$sortObj = $ws1.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws1.Range("D2:D10"))
$sortObj.SetRange($ws1.Range("A2:E10"))
$sortObj.Header = -4142
$sortObj.Apply()

# Rebuild hyperlinks on the e-mail column for everyone except Jensen.
$ws1.Hyperlinks.Add($ws1.Range("C2"), "mailto:amorales@global.com")
$ws1.Range("C2").Style = "Hyperlink"
$ws1.Hyperlinks.Add($ws1.Range("C3"), "mailto:vsharma@global.com")
$ws1.Range("C3").Style = "Hyperlink"
$ws1.Hyperlinks.Add($ws1.Range("C4"), "mailto:dwilliams@global.com")
$ws1.Range("C4").Style = "Hyperlink"
$ws1.Hyperlinks.Add($ws1.Range("C5"), "mailto:epetrov@global.com")
$ws1.Range("C5").Style = "Hyperlink"
$ws1.Hyperlinks.Add($ws1.Range("C7"), "mailto:ahassan@global.com")
$ws1.Range("C7").Style = "Hyperlink"

# Jensen (row 6) and the now-empty rows 8:10 keep the Hyperlink look on C
# without an actual link, same as before the edit.
$ws1.Range("C6").Style = "Hyperlink"
$ws1.Range("C8").Style = "Hyperlink"
$ws1.Range("C9").Style = "Hyperlink"
$ws1.Range("C10").Style = "Hyperlink"

# Explicit widths for the Email Address / Title columns.
$ws1.Columns.Item(3).ColumnWidth = 21.998697916666668
$ws1.Columns.Item(4).ColumnWidth = 17.830729166666668

# ---------------------------------------------------------------------
# Departments sheet -- drop the duplicate "Human Resources" row
# ---------------------------------------------------------------------
$ws2.Rows.Item(5).Delete()

# ---------------------------------------------------------------------
# View state -- Employees tab is now the active/selected one
# ---------------------------------------------------------------------
$ws2.Range("D16").Select()
$ws1.Activate()
$ws1.Range("A10").Select()
